# "replico imagen error y segunda diapo para commit"
#
# Inserts a new slide ("Diapo de segundo commit") as the 3rd slide of the
# deck (i.e. right after the "DIAPO PROFE" slide and before "¿Qué es Git?").
# The new slide uses the same "Title and Content" layout (layout index 2,
# matching slideLayout2.xml / "Título y objetos") used by the other content
# slides, has a title placeholder with the text "Diapo de segundo commit"
# and an empty content placeholder.

$p = $ppt.ActivePresentation

# Insert the new slide at position 3 using the Title-and-Content layout.
$s = $p.Slides.Add(3, 2)

# Fill in the title placeholder text; leave the content placeholder empty.
$s.Shapes.Title.TextFrame.TextRange.Text = "Diapo de segundo commit"
